$d = $word.ActiveDocument

# The "Ken Alleyne" section currently has an empty body paragraph (just a
# tab character). Find that heading, then grab the content paragraph that
# follows it, and append this week's update text right after the existing
# tab run (but before the paragraph mark).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Ken Alleyne") {
        $target = $p.Next()
        break
    }
}

$insertAt = $target.Range.End - 1

$chunks = @(
    "This week I ",
    "briefly looked into an issue we were having with ",
    "displaying user data o",
    "n the front end, and ",
    "offered assistance with",
    " implementing the functionality that ",
    "allows users to update their own tasks."
)

foreach ($chunk in $chunks) {
    $ip = $d.Range($insertAt, $insertAt)
    $ip.InsertAfter($chunk)
    $insertAt = $insertAt + $chunk.Length
}
